$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the data set. This shifts the
# existing rows 89-123 down to 90-124 and inserts the new record at row 89.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly record.
$ws.Range("A89").Value = 7
$ws.Range("B89").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C89").Value = "Ñuble"
$ws.Range("D89").Value = 44985
$ws.Range("E89").Value = 16
$ws.Range("F89").Value = 100112031
$ws.Range("G89").Value = "Poroto verde"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 60
$ws.Range("K89").Value = 26000
$ws.Range("L89").Value = 27000
$ws.Range("M89").Value = 26500
$ws.Range("N89").Value = "$/saco 25 kilos"
$ws.Range("O89").Value = "Provincia de Diguillín"
$ws.Range("P89").Value = 1060
$ws.Range("Q89").Value = 25
$ws.Range("R89").Value = "Hortaliza"
